$wb = $excel.ActiveWorkbook

$wsDeporte = $wb.Worksheets.Item("Deporte")
$wsHistoria = $wb.Worksheets.Item("Historia")

# Update the Historia questions to the new wording (wrapped with inverted question marks).
# The shared-strings table is compacted/reindexed automatically by the runtime,
# which also naturally reindexes the color values used on every sheet's column E.
$wsHistoria.Range("C1").Value = "¿Donde nació Jesús?"
$wsHistoria.Range("C2").Value = "¿Que representan los colores de la bandera Argentina?"
$wsHistoria.Range("C3").Value = "¿Donde falleció San Martín?"

# Widen column C on the Historia sheet to fit the longer question text.
$wsHistoria.Columns.Item(3).ColumnWidth = 70.15

# Move the selection on the Historia sheet (no longer the active tab) to C5.
$wsHistoria.Range("C5").Select()

# Make "Deporte" the active/selected tab, keeping its existing selection (E3).
$wsDeporte.Activate()
